$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.141.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5212"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06326"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07539"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.666.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.436"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5501"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.145.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.749"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.223"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1239"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.471"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06305"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.348"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.525"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.647"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.005"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6036"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.401"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.757"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.113.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01614"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8632"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.055"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4242"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.83%  "
